# Update countries & provincias Spain
# - Refresh the "datos actualizados" timestamp (11:35 -> 12:05)
# - New daily figures for Rumania/Austria/Japon (rows 37-39): Rumania's
#   updated numbers land on row 37, pushing Austria's old row-37 numbers to
#   row 38 and Japon's old row-38 numbers to row 39.
# - Updated daily figures for Malta (row 124) and Etiopia (row 141)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header timestamp
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 13 de Mayo de 2020 a las 12:05"

# Row 37 -> Rumania's new numbers
$ws.Cells.Item(37, 1).Value = "Rumania"
$ws.Cells.Item(37, 2).Value = 16002
$ws.Cells.Item(37, 3).Value = 224
$ws.Cells.Item(37, 4).Value = 7961
$ws.Cells.Item(37, 5).Value = 7025
$ws.Cells.Item(37, 6).Value = 228
$ws.Cells.Item(37, 7).Value = 14
$ws.Cells.Item(37, 8).Value = 1016

# Row 38 -> Austria (takes over what used to be row 37's numbers)
$ws.Cells.Item(38, 1).Value = "Austria"
$ws.Cells.Item(38, 2).Value = 15997
$ws.Cells.Item(38, 3).Value = 36
$ws.Cells.Item(38, 4).Value = 14304
$ws.Cells.Item(38, 5).Value = 1069
$ws.Cells.Item(38, 6).Value = 55
$ws.Cells.Item(38, 7).Value = 1
$ws.Cells.Item(38, 8).Value = 624

# Row 39 -> Japon (takes over what used to be row 38's numbers)
$ws.Cells.Item(39, 1).Value = "Japon"
$ws.Cells.Item(39, 2).Value = 15968
$ws.Cells.Item(39, 3).Value = 0
$ws.Cells.Item(39, 4).Value = 8531
$ws.Cells.Item(39, 5).Value = 6780
$ws.Cells.Item(39, 6).Value = 249
$ws.Cells.Item(39, 7).Value = 0
$ws.Cells.Item(39, 8).Value = 657

# Row 124 -> Malta updated figures
$ws.Cells.Item(124, 5).Value = 66
$ws.Cells.Item(124, 7).Value = 1
$ws.Cells.Item(124, 8).Value = 6

# Row 141 -> Etiopia updated figures
$ws.Cells.Item(141, 2).Value = 263
$ws.Cells.Item(141, 3).Value = 2
$ws.Cells.Item(141, 4).Value = 108
